# Adds the "2022-Q3" quarter to the BXP workbook:
#  1. Insert a new "2022-Q3" worksheet right after the "总计" summary sheet,
#     copying the layout/formatting of the existing "2022-Q1" sheet and
#     filling in the new quarter's figures.
#  2. Insert a new row at the top of the "总计" summary sheet for 2022-Q3
#     (existing quarters shift down one row, keeping their own values).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Item("2022-Q1")

# --- 1. Create the new "2022-Q3" sheet right before "2022-Q1" -------------
$q3 = $wb.Worksheets.Add($q1)
$q3.Name = "2022-Q3"

# Bring over headers/layout/formatting from the "2022-Q1" sheet.
$q1.Range("A1:H2").Copy($q3.Range("A1"))

# Fill in the 2022-Q3 figures.
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "320017"
$q3.Cells.Item(2, 3).Value = "诺安全球收益不动产（QDII）"
$q3.Cells.Item(2, 4).Value = "0.23"
$q3.Cells.Item(2, 5).Value = "73.76"
$q3.Cells.Item(2, 6).Value = "4.66"
$q3.Cells.Item(2, 7).Value = "0.0107"
$q3.Cells.Item(2, 8).Value = 9

# --- 2. Shift the "总计" summary rows down and add the 2022-Q3 entry ------
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 4).Value = 0.01

$summary.Cells.Item(3, 2).Value = "2022-Q1"
$summary.Cells.Item(3, 4).Value = 0.02

$summary.Cells.Item(4, 2).Value = "2021-Q4"
$summary.Cells.Item(4, 4).Value = 0.03

$summary.Cells.Item(5, 2).Value = "2021-Q3"
$summary.Cells.Item(5, 4).Value = 0.02

$summary.Cells.Item(6, 2).Value = "2021-Q2"
$summary.Cells.Item(6, 4).Value = 0.03

$summary.Cells.Item(7, 2).Value = "2021-Q1"
$summary.Cells.Item(7, 4).Value = 0.03

# New eighth row, matching the formatting of the row above it.
$summary.Cells.Item(7, 1).Copy()
$summary.Cells.Item(8, 1).PasteSpecial(-4122)
$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(8, 2).Value = "2020-Q4"
$summary.Cells.Item(8, 3).Value = 1
$summary.Cells.Item(8, 4).Value = 0.03
